$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row for 2022-Q4 at the top of the data,
#    shifting the existing 2022-Q3 and 2022-Q2 rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift existing data down first (bottom-up so values aren't clobbered
# before they are read).
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 8
$total.Cells.Item(4,4).Value = 0.83

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 3
$total.Cells.Item(3,4).Value = 0.05

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 15
$total.Cells.Item(2,4).Value = 0.85

# Column A carries the bordered/centred style already used on rows 2 & 3;
# copy it onto the newly-created row 4 as well.
$total.Range("A2").Copy()
$total.Range("A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right after "总计" (so the tab order
#    becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2).
# ---------------------------------------------------------------------------
$newQ4 = $wb.Worksheets.Add($null, $total, 1, $null)
$newQ4.Name = "2022-Q4"

# Fetch the "2022-Q3" sheet reference AFTER the insert: worksheet handles in
# this engine resolve by tab position, and the insert above shifted "2022-Q3"
# from slot 2 to slot 3.
$q3 = $wb.Worksheets.Item(3)   # existing "2022-Q3" sheet - used as a style donor

# Copy header-row and column-A formatting from the 2022-Q3 sheet so the new
# sheet matches the same look (bold, centred, thin border = style "s=2").
$q3.Range("B1:H1").Copy()
$newQ4.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2").Copy()
$newQ4.Range("A2:A16").PasteSpecial(-4122)

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $newQ4.Cells.Item(1, $col + 2).Value = $headers[$col]
}

$q4Data = @(
    ,@("014915","财通匠心优选一年持有期混合A","5.48","91.42","6.24","0.3420",5)
    ,@("009062","财通智慧成长混合A","2.05","86.49","4.54","0.0931",9)
    ,@("009063","财通智慧成长混合C","1.74","86.49","4.54","0.0790",9)
    ,@("202019","南方策略优化混合","2.80","93.64","2.75","0.0770",2)
    ,@("006693","金信消费升级股票C","0.97","87.56","6.07","0.0589",7)
    ,@("166109","信澳量化先锋混合（LOF）A","0.90","94.26","5.80","0.0522",2)
    ,@("006692","金信消费升级股票A","0.75","87.56","6.07","0.0455",7)
    ,@("014916","财通匠心优选一年持有期混合C","0.59","91.42","6.24","0.0368",5)
    ,@("002862","金信量化精选灵活配置混合","0.46","80.41","4.13","0.0190",7)
    ,@("001244","华泰柏瑞量化智慧灵活配置混合A","2.89","93.57","0.65","0.0188",5)
    ,@("000757","华富智慧城市灵活配置混合","0.48","92.86","3.37","0.0162",5)
    ,@("166110","信澳量化先锋混合（LOF）C","0.20","94.26","5.80","0.0116",2)
    ,@("006104","华泰柏瑞量化智慧灵活配置混合C","0.38","93.57","0.65","0.0025",5)
    ,@("006857","蜂巢卓睿灵活配置混合A","0.07","68.96","1.35","0.0009",1)
    ,@("006858","蜂巢卓睿灵活配置混合C","0.03","68.96","1.35","0.0004",1)
)

# Columns B-G hold strings that look numeric ("014915", "5.48", "0.3420", …).
# Force text entry so leading zeros / trailing zeros survive verbatim,
# matching the source file's inlineStr cells.
$textRange = $newQ4.Range("B2:G16")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = $i + 2
    $row = $q4Data[$i]
    $newQ4.Cells.Item($r, 1).Value = $i
    $newQ4.Cells.Item($r, 2).Value = $row[0]
    $newQ4.Cells.Item($r, 3).Value = $row[1]
    $newQ4.Cells.Item($r, 4).Value = $row[2]
    $newQ4.Cells.Item($r, 5).Value = $row[3]
    $newQ4.Cells.Item($r, 6).Value = $row[4]
    $newQ4.Cells.Item($r, 7).Value = $row[5]
    $newQ4.Cells.Item($r, 8).Value = $row[6]
}

# Drop the text-coercion format back to the default (unstyled) cell style so
# these cells end up styleless, just like the source file.
$textRange.Style = "Normal"

Write-Output "edit applied"
